# Reverse the comma-separated "Recorded By" list in column G whenever it
# starts with the literal token "System" (e.g. "System, foo@bar.com"
# becomes "foo@bar.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($null -ne $val -and $val -is [string] -and $val.Length -gt 0) {
        $parts = $val -split ', '
        if ($parts.Count -gt 1 -and $parts[0] -eq 'System') {
            $n = $parts.Count
            $reversed = @()
            for ($i = $n - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }
            $cell.Value = [string]::Join(', ', $reversed)
        }
    }
}
